$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as the new first data row of this
# weekly block (row 143), pushing all the existing rows from 143..195 down
# to 144..196.
$ws.Rows.Item(143).Insert()

$ws.Cells.Item(143, 1).Value  = 3
$ws.Cells.Item(143, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(143, 3).Value  = "Coquimbo"
$ws.Cells.Item(143, 4).Value  = 44468
$ws.Cells.Item(143, 5).Value  = 5
$ws.Cells.Item(143, 6).Value  = 100114013
$ws.Cells.Item(143, 7).Value  = "Zanahoria"
$ws.Cells.Item(143, 8).Value  = "Sin especificar"
$ws.Cells.Item(143, 9).Value  = "Primera"
$ws.Cells.Item(143, 10).Value = 250
$ws.Cells.Item(143, 11).Value = 8000
$ws.Cells.Item(143, 12).Value = 8000
$ws.Cells.Item(143, 13).Value = 8000
$ws.Cells.Item(143, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(143, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(143, 16).Value = 400
$ws.Cells.Item(143, 17).Value = 20
$ws.Cells.Item(143, 18).Value = "Hortaliza"
